$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure cells are treated as text so formatted numeric-looking strings
# (e.g. "1.00", "0.840", "62.538.77") are preserved verbatim, matching
# the source data which stores these as inline text, not numbers.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.538.77"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -2.87%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.380.64"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -3.68%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "572.66"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -3.26%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "125.49"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -6.82%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.380.44"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -3.65%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.475"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -2.52%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -4.22%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -4.56%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.375"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -3.34%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.955.10"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -3.72%  "
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.92%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.377.88"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -3.77%  "
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -6.04%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.562.75"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -2.82%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "24.40"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -5.24%  "
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -7.34%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.62"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -2.31%  "
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -4.15%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "372.76"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -5.19%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.514.58"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -3.70%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "71.62"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -4.01%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000106"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -10.17%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.17%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -6.77%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -7.18%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.77"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -5.97%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.01%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -6.66%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.410.18"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.148"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -5.92%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.29"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -1.51%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "166.67"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -0.15%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.63"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -4.64%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -5.53%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.00"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.06%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -5.82%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "41.49"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -1.25%  "
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -4.92%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -8.85%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -7.29%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.08"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -8.33%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.59"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -3.08%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.239.72"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -5.67%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.840"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -7.72%  "
